$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking snapshot refresh: update Price (D) / Volume(1h) (E) for each row,
# and for two rank swaps, the Coin (B) / Link (C) cells as well.
# Price values that look like plain numbers are entered with a leading
# apostrophe so Excel stores them as literal text (matching the sheet's
# existing inlineStr convention) instead of auto-converting to a number.

# Row 2
$ws.Range("D2").Value = "72.183.48"
$ws.Range("E2").Value = "  +0.84%  "

# Row 3
$ws.Range("D3").Value = "2.710.89"
$ws.Range("E3").Value = "  +3.25%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").Value = "'599.74"
$ws.Range("E5").Value = "  -1.11%  "

# Row 6
$ws.Range("D6").Value = "'176.04"
$ws.Range("E6").Value = "  -1.71%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").Value = "'0.525"
$ws.Range("E8").Value = "  -0.18%  "

# Row 9
$ws.Range("D9").Value = "2.711.04"
$ws.Range("E9").Value = "  +3.28%  "

# Row 10
$ws.Range("D10").Value = "'0.169"
$ws.Range("E10").Value = "  +0.76%  "

# Row 11
$ws.Range("E11").Value = "  +2.62%  "

# Row 12
$ws.Range("D12").Value = "'0.355"
$ws.Range("E12").Value = "  +2.01%  "

# Row 13
$ws.Range("D13").Value = "'5.02"
$ws.Range("E13").Value = "  -0.45%  "

# Row 14
$ws.Range("D14").Value = "3.208.52"
$ws.Range("E14").Value = "  +2.36%  "

# Row 15
$ws.Range("E15").Value = "  -0.20%  "

# Row 16
$ws.Range("D16").Value = "71.978.12"
$ws.Range("E16").Value = "  +0.57%  "

# Row 17
$ws.Range("D17").Value = "'26.35"
$ws.Range("E17").Value = "  -0.69%  "

# Row 18
$ws.Range("D18").Value = "2.711.13"
$ws.Range("E18").Value = "  +3.12%  "

# Row 19
$ws.Range("D19").Value = "'12.28"
$ws.Range("E19").Value = "  +7.34%  "

# Row 20
$ws.Range("D20").Value = "'8.13"
$ws.Range("E20").Value = "  +2.13%  "

# Row 21
$ws.Range("D21").Value = "'373.59"
$ws.Range("E21").Value = "  -2.39%  "

# Row 22
$ws.Range("D22").Value = "'4.17"
$ws.Range("E22").Value = "  +1.25%  "

# Row 23
$ws.Range("E23").Value = "  +2.57%  "

# Row 24
$ws.Range("D24").Value = "'72.43"
$ws.Range("E24").Value = "  -0.16%  "

# Row 25
$ws.Range("E25").Value = "  -0.03%  "

# Row 26
$ws.Range("D26").Value = "'4.37"
$ws.Range("E26").Value = "  -1.47%  "

# Row 27
$ws.Range("D27").Value = "'9.87"
$ws.Range("E27").Value = "  -0.98%  "

# Row 28
$ws.Range("D28").Value = "2.851.52"
$ws.Range("E28").Value = "  +3.20%  "

# Row 29
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.18%  "

# Row 30
$ws.Range("D30").Value = "0.0₂01000"
$ws.Range("E30").Value = "  +4.19%  "

# Row 31
$ws.Range("D31").Value = "'8.14"
$ws.Range("E31").Value = "  +1.32%  "

# Row 32
$ws.Range("D32").Value = "'508.45"
$ws.Range("E32").Value = "  -6.36%  "

# Row 33
$ws.Range("E33").Value = "  -1.62%  "

# Row 34
$ws.Range("E34").Value = "  -0.06%  "

# Row 35
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.01%  "

# Row 36
$ws.Range("D36").Value = "'164.11"
$ws.Range("E36").Value = "  -1.26%  "

# Row 37
$ws.Range("D37").Value = "'19.72"
$ws.Range("E37").Value = "  +2.72%  "

# Row 38
$ws.Range("D38").Value = "'19.11"
$ws.Range("E38").Value = "  -0.27%  "

# Row 39
$ws.Range("E39").Value = "  -0.36%  "

# Row 40
$ws.Range("E40").Value = "  -4.50%  "

# Row 41
$ws.Range("D41").Value = "'1.81"
$ws.Range("E41").Value = "  -2.76%  "

# Row 42
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.03%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").Value = "'5.07"
$ws.Range("E43").Value = "  +0.74%  "

# Row 44
$ws.Range("D44").Value = "'2.58"
$ws.Range("E44").Value = "  -1.96%  "

# Row 45
$ws.Range("E45").Value = "  +0.77%  "

# Row 46
$ws.Range("D46").Value = "'157.38"
$ws.Range("E46").Value = "  +4.31%  "

# Row 47
$ws.Range("D47").Value = "'39.46"
$ws.Range("E47").Value = "  +0.68%  "

# Row 48
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'0.563"
$ws.Range("E48").Value = "  +5.25%  "

# Row 49
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").Value = "'3.75"
$ws.Range("E49").Value = "  +3.16%  "

# Row 50
$ws.Range("E50").Value = "  +5.31%  "

# Row 51
$ws.Range("D51").Value = "'0.0767"
$ws.Range("E51").Value = "  +1.50%  "
